$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1280.079
$ws.Range("D2").Value = 334.028
$ws.Range("C3").Value = 3397.134
$ws.Range("D3").Value = 448.809
$ws.Range("C4").Value = 2313.455
$ws.Range("D4").Value = 337.8939999999999
$ws.Range("C11").Value = 861.854
$ws.Range("C12").Value = 1407.019
$ws.Range("C13").Value = 1119.255
$ws.Range("C14").Value = 2028.106
$ws.Range("C16").Value = 2826.968
$ws.Range("C23").Value = 1206.8
$ws.Range("C24").Value = 2332.294
$ws.Range("C25").Value = 1725.248
$ws.Range("C26").Value = 263.7859999999999
$ws.Range("C29").Value = 724.1210000000001
$ws.Range("C30").Value = 976.073
$ws.Range("C31").Value = 772.2090000000001
$ws.Range("C35").Value = 479.79
$ws.Range("C36").Value = 916.442
$ws.Range("C37").Value = 686.8950000000001
$ws.Range("C50").Value = 762.275
$ws.Range("C51").Value = 610.885
$ws.Range("C52").Value = 1436.499
$ws.Range("C53").Value = 1183.668
$ws.Range("C54").Value = 242.098
$ws.Range("D54").Value = 137.731
$ws.Range("C55").Value = 645.9359999999999
$ws.Range("D55").Value = 363.638
$ws.Range("C56").Value = 383.966
$ws.Range("C63").Value = 861.3300000000003
$ws.Range("D63").Value = 586.0069999999999
$ws.Range("C64").Value = 1755.687
$ws.Range("C65").Value = 1278.732
$ws.Range("C66").Value = 1232.941
$ws.Range("C67").Value = 2404.332
$ws.Range("C68").Value = 1740.853
$ws.Range("C69").Value = 702.2090000000001
$ws.Range("D69").Value = 200.2130000000001
$ws.Range("C70").Value = 1341.706
$ws.Range("D70").Value = 250.376
$ws.Range("C71").Value = 991.724
$ws.Range("D71").Value = 221.407
$ws.Range("C72").Value = 183.495
$ws.Range("D72").Value = 154.174
$ws.Range("C73").Value = 278.952
$ws.Range("D73").Value = 142.06
$ws.Range("C74").Value = 222.38
$ws.Range("D74").Value = 142.155
$ws.Range("C81").Value = 1520.329
$ws.Range("C82").Value = 3214.981
$ws.Range("C83").Value = 2353.129
$ws.Range("C87").Value = 1170.017
$ws.Range("C88").Value = 2094.6
$ws.Range("C89").Value = 1555.849
$ws.Range("C90").Value = 997.1990000000001
$ws.Range("C91").Value = 1837.148
$ws.Range("C92").Value = 1364.685
$ws.Range("C93").Value = 1015.775
$ws.Range("D93").Value = 176.354
$ws.Range("C94").Value = 1698.117
$ws.Range("C95").Value = 1355.724
$ws.Range("D95").Value = 204.065
